# CompStat weekly update: new crime data collected.
# Updates the report header (volume/date) text and the precinct crime-stat
# table (rows 14-33) to the new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text: bump the report's issue number and the covered week dates.
# (These source cells are rich-text shared strings whose runs all share
# the same font as the cell's own style, so a plain value update keeps
# the visible formatting identical.)
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  24"
$ws.Range("C9").Value = "Report Covering the Week  6/10/2024  Through  6/16/2024"

# ---------------------------------------------------------------------
# Cells whose content TYPE changes (number <-> the "0"/"***.*" text
# placeholders used for n/a figures). Excel applies a different cell
# style to text vs. numeric entries in this table (General vs. numeric
# formats), so we copy an already-correctly-styled neighbor cell into
# place (which carries both the right style AND, where it already
# matches, the right value), then overwrite the value when it still
# needs to differ from the donor.
# ---------------------------------------------------------------------

# --- become the "0" text placeholder (donor: C22, already "0") ---
$ws.Range("C22").Copy($ws.Range("C14"))
$ws.Range("C22").Copy($ws.Range("C28"))
$ws.Range("C22").Copy($ws.Range("D28"))
$ws.Range("C22").Copy($ws.Range("F31"))
$ws.Range("C22").Copy($ws.Range("D33"))

# --- become the "***.*" text placeholder (donor: E22, already "***.*") ---
$ws.Range("E22").Copy($ws.Range("E28"))
$ws.Range("E22").Copy($ws.Range("E33"))

# --- become numeric (donor: F14, numeric style), then set the real value ---
$ws.Range("F14").Copy($ws.Range("C15"))
$ws.Range("C15").Value = 2

$ws.Range("F14").Copy($ws.Range("C27"))
$ws.Range("C27").Value = 2

# ---------------------------------------------------------------------
# Remaining value-only updates (cell style/number format unchanged).
# ---------------------------------------------------------------------
$ws.Range("F15").Value = 3
$ws.Range("I15").Value = 10
$ws.Range("K15").Value = 11.111111111111
$ws.Range("L15").Value = -9.090909090909
$ws.Range("M15").Value = -23.076923076923
$ws.Range("N15").Value = -41.176470588235
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -60
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 94
$ws.Range("J16").Value = 91
$ws.Range("K16").Value = 3.296703296703
$ws.Range("L16").Value = 6.818181818181
$ws.Range("M16").Value = -12.962962962963
$ws.Range("N16").Value = -75.835475578406
$ws.Range("D17").Value = 11
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 35
$ws.Range("G17").Value = 38
$ws.Range("H17").Value = -7.894736842105
$ws.Range("I17").Value = 172
$ws.Range("J17").Value = 153
$ws.Range("K17").Value = 12.418300653594
$ws.Range("L17").Value = 29.323308270676
$ws.Range("M17").Value = 160.606060606061
$ws.Range("N17").Value = 4.878048780487
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = -41.666666666666
$ws.Range("I18").Value = 54
$ws.Range("J18").Value = 50
$ws.Range("K18").Value = 8
$ws.Range("L18").Value = -15.625
$ws.Range("M18").Value = -58.139534883720
$ws.Range("N18").Value = -89.868667917448
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = -68.75
$ws.Range("F19").Value = 32
$ws.Range("G19").Value = 54
$ws.Range("H19").Value = -40.740740740740
$ws.Range("I19").Value = 228
$ws.Range("J19").Value = 289
$ws.Range("K19").Value = -21.107266435986
$ws.Range("L19").Value = -29.629629629629
$ws.Range("M19").Value = 53.020134228187
$ws.Range("N19").Value = -11.284046692607
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 12
$ws.Range("E20").Value = -75
$ws.Range("F20").Value = 27
$ws.Range("G20").Value = 28
$ws.Range("H20").Value = -3.571428571428
$ws.Range("I20").Value = 149
$ws.Range("J20").Value = 109
$ws.Range("K20").Value = 36.697247706422
$ws.Range("L20").Value = 29.565217391304
$ws.Range("M20").Value = 17.322834645669
$ws.Range("N20").Value = -90.274151436031
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 48
$ws.Range("E21").Value = -50
$ws.Range("F21").Value = 120
$ws.Range("G21").Value = 147
$ws.Range("H21").Value = -18.367346938775
$ws.Range("I21").Value = 710
$ws.Range("J21").Value = 702
$ws.Range("K21").Value = 1.139601139601
$ws.Range("L21").Value = -3.532608695652
$ws.Range("M21").Value = 19.730185497470
$ws.Range("N21").Value = -75.567790777701
$ws.Range("F22").Value = 1
$ws.Range("L22").Value = 50
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 34
$ws.Range("E24").Value = -44.117647058823
$ws.Range("F24").Value = 81
$ws.Range("G24").Value = 112
$ws.Range("H24").Value = -27.678571428571
$ws.Range("I24").Value = 545
$ws.Range("J24").Value = 593
$ws.Range("K24").Value = -8.094435075885
$ws.Range("L24").Value = -15.109034267912
$ws.Range("M24").Value = 89.236111111111
$ws.Range("C25").Value = 15
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = 36.363636363636
$ws.Range("F25").Value = 45
$ws.Range("G25").Value = 45
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 274
$ws.Range("J25").Value = 241
$ws.Range("K25").Value = 13.692946058091
$ws.Range("L25").Value = -4.529616724738
$ws.Range("C26").Value = 9
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = -10
$ws.Range("F26").Value = 50
$ws.Range("G26").Value = 48
$ws.Range("H26").Value = 4.166666666666
$ws.Range("I26").Value = 260
$ws.Range("J26").Value = 251
$ws.Range("K26").Value = 3.585657370517
$ws.Range("L26").Value = 25.603864734299
$ws.Range("M26").Value = 15.044247787610
$ws.Range("F27").Value = 3
$ws.Range("I27").Value = 14
$ws.Range("K27").Value = 7.692307692307
$ws.Range("L27").Value = -26.315789473684
$ws.Range("F28").Value = 3
$ws.Range("H28").Value = 0
$ws.Range("L28").Value = -25.925925925925
$ws.Range("L29").Value = -66.666666666666
$ws.Range("L30").Value = -62.5
$ws.Range("G33").Value = 2
